$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $origStyle = $Cell.Style
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "320.53"
Set-TextValue $ws.Range("E2") "-3.58%"
Set-TextValue $ws.Range("D3") "42.89"
Set-TextValue $ws.Range("E3") "-6.79%"
Set-TextValue $ws.Range("D4") "5.187"
Set-TextValue $ws.Range("E4") "-8.62%"
Set-TextValue $ws.Range("D5") "0.08144"
Set-TextValue $ws.Range("E5") "-2.78%"
Set-TextValue $ws.Range("D6") "4.342"
Set-TextValue $ws.Range("E6") "-3.10%"
Set-TextValue $ws.Range("D7") "1.790"
Set-TextValue $ws.Range("E7") "-12.35%"
Set-TextValue $ws.Range("D8") "0.9523"
Set-TextValue $ws.Range("E8") "-3.88%"
Set-TextValue $ws.Range("D9") "0.1115"
Set-TextValue $ws.Range("E9") "-3.41%"
Set-TextValue $ws.Range("D10") "0.1845"
Set-TextValue $ws.Range("E10") "-4.99%"
Set-TextValue $ws.Range("D11") "0.04668"
Set-TextValue $ws.Range("E11") "-0.11%"
Set-TextValue $ws.Range("D12") "0.09341"
Set-TextValue $ws.Range("E12") "-6.55%"
Set-TextValue $ws.Range("D13") "7.451"
Set-TextValue $ws.Range("E13") "-28.31%"
Set-TextValue $ws.Range("D14") "0.1061"
Set-TextValue $ws.Range("E14") "0.17%"
Set-TextValue $ws.Range("D15") "0.001289"
Set-TextValue $ws.Range("E15") "0.95%"
Set-TextValue $ws.Range("D16") "0.005916"
Set-TextValue $ws.Range("E16") "-2.03%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D17") "3.365"
Set-TextValue $ws.Range("E17") "-0.17%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D18") "2.530"
Set-TextValue $ws.Range("E18") "-1.75%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue $ws.Range("D19") "0.3364"
Set-TextValue $ws.Range("E19") "-0.04%"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue $ws.Range("D20") "0.1394"
Set-TextValue $ws.Range("E20") "-0.52%"
$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextValue $ws.Range("D21") "0.2628"
Set-TextValue $ws.Range("E21") "-0.82%"
$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws.Range("D22") "0.04188"
Set-TextValue $ws.Range("E22") "-0.73%"
$ws.Range("B23").Value = "BitKan"
$ws.Range("C23").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue $ws.Range("D23") "0.001253"
Set-TextValue $ws.Range("E23") "-4.24%"
$ws.Range("B24").Value = "HotbitToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue $ws.Range("D24") "0.004325"
Set-TextValue $ws.Range("E24") "-6.77%"
Set-TextValue $ws.Range("D25") "0.0001111"
Set-TextValue $ws.Range("E25") "-13.31%"
Set-TextValue $ws.Range("D26") "0.0002992"
Set-TextValue $ws.Range("E26") "-20.12%"
Set-TextValue $ws.Range("D38") "0.02583"
Set-TextValue $ws.Range("E38") "-7.90%"
Set-TextValue $ws.Range("D39") "0.05532"
Set-TextValue $ws.Range("E39") "-4.44%"
Set-TextValue $ws.Range("D40") "0.007825"
Set-TextValue $ws.Range("E40") "0.93%"
Set-TextValue $ws.Range("D41") "0.1392"
Set-TextValue $ws.Range("E41") "-3.15%"
Set-TextValue $ws.Range("D42") "0.006621"
Set-TextValue $ws.Range("E42") "-9.20%"
Set-TextValue $ws.Range("D43") "0.002118"
Set-TextValue $ws.Range("E43") "7.21%"
Set-TextValue $ws.Range("D44") "0.008480"
Set-TextValue $ws.Range("E44") "-6.23%"
Set-TextValue $ws.Range("D45") "0.3460"
Set-TextValue $ws.Range("E45") "1.37%"
Set-TextValue $ws.Range("D46") "0.00006991"
Set-TextValue $ws.Range("E46") "-5.33%"
Set-TextValue $ws.Range("E47") "0.28%"
Set-TextValue $ws.Range("D48") "0.003478"
Set-TextValue $ws.Range("E48") "-0.78%"
Set-TextValue $ws.Range("D49") "0.003545"
Set-TextValue $ws.Range("E49") "1.20%"
Set-TextValue $ws.Range("D50") "0.00002109"
Set-TextValue $ws.Range("E50") "0.28%"
Set-TextValue $ws.Range("D51") "0.0002009"
Set-TextValue $ws.Range("E51") "0.28%"
